# Moved import_upload file to shrine
#
# The "Investor" column (E) — header + per-row investor-name values — is no
# longer populated by the importer, so clear its contents (header "Investor",
# "Kalaari Capital", "Accel") while leaving the column's formatting/style in
# place. Also update the active selection to match the column that was
# touched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Investor" column's header + data (row 1 header, rows 2-3 values)
# but keep the cells (and their style) in place.
$ws.Range("E1:E3").ClearContents()

# Reflect the edited column in the active selection (was a single cell E4).
$ws.Range("E:E").Select()
